# Swap the contents of columns B:AC between specific row pairs on the active sheet.
# Column A (the leading index column) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 2   # Column B
$lastCol  = 29  # Column AC

$rowPairs = @(
    @(29, 30),
    @(36, 37),
    @(99, 100),
    @(111, 112)
)

foreach ($pair in $rowPairs) {
    $row1 = $pair[0]
    $row2 = $pair[1]

    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cell1 = $ws.Cells.Item($row1, $col)
        $cell2 = $ws.Cells.Item($row2, $col)

        $v1 = $cell1.Value()
        $v2 = $cell2.Value()

        $cell1.Value = $v2
        $cell2.Value = $v1
    }
}
